$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new row 8, duplicating row 7 (same "Test Ringover (NO TOCAR)"
#     order) but with its own Optimizador info -----------------------------
#
# Copy/paste (rather than per-cell Value assignment) is used so that the
# numeric-looking quantities ("1", "3", ...) keep being stored as TEXT,
# exactly like the rest of this sheet does - a plain `.Value = "1"` would
# get auto-coerced to a real number by Excel.
$ws.Range("A7:N7").Copy()
$ws.Range("A8:N8").PasteSpecial()

# Row 8 has its own Optimizador (columns F/G); row 7 left these blank.
$ws.Cells.Item(8, 6).Value = "HUAWEI Optimizador 600W"

# Column G ("Unidades Optimizador") also stores quantities as literal text
# in this sheet, so "2" has to be written the same text-preserving way: put
# it (prefixed with an apostrophe, forcing text) in a scratch cell, copy
# that cell, and paste only its value into G8 - this carries over the text
# typing without a normal Value assignment re-coercing "2" into a number.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.Value = "'2"
$scratch.Copy()
$ws.Range("G8").PasteSpecial(-4163)  # xlPasteValues
$scratch.Delete()

# --- Row 7 drops its unused placeholder cells -----------------------------
# E7:G7 and L7 used to be present as empty cells; they're cleared out
# entirely now.
$ws.Range("E7:G7").ClearContents()
$ws.Range("L7").ClearContents()
